# Update the MEC-2B schedule sheet so that the assignments recorded in
# columns B, C and E for rows 2, 3, 4, 6, 7 and 8 reflect the refreshed
# schedule ("Só para o Túlio pegar atualizado").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "-"
$ws.Range("E2").Value = "[Elcio Dec.-Des. Maq. Cad._T1-2B, Elcio Dec.-Des. Maq. Cad._T1-2B]"

$ws.Range("B3").Value = "Maria Celeste-Máquinas Térmicas e de Fluxo"
$ws.Range("C3").Value = "-"
$ws.Range("E3").Value = "[Carlos-Tornearia-2B, Victor-Ajustagem-2B, Elaine-Metalografia-2B, Emerson-Elet. Dig. Bas.-2B]"

$ws.Range("B4").Value = "Maria Celeste-Máquinas Térmicas e de Fluxo"
$ws.Range("C4").Value = "-"
$ws.Range("E4").Value = "[Carlos-Tornearia-2B, Victor-Ajustagem-2B, Elaine-Metalografia-2B, Emerson-Elet. Dig. Bas.-2B]"

$ws.Range("B6").Value = "Gilberto-Mec. Tec. Res. Mat."
$ws.Range("C6").Value = "-"
$ws.Range("E6").Value = "[Carlos-Tornearia-2B, Victor-Ajustagem-2B, Elaine-Metalografia-2B, Emerson-Elet. Dig. Bas.-2B]"

$ws.Range("B7").Value = "[Elcio Dec.-Des. Maq. Cad._T2-2B, Elcio Dec.-Des. Maq. Cad._T2-2B]"
$ws.Range("C7").Value = "-"
$ws.Range("E7").Value = "[Carlos-Tornearia-2B, Victor-Ajustagem-2B, Elaine-Metalografia-2B, Emerson-Elet. Dig. Bas.-2B]"

$ws.Range("C8").Value = "-"
$ws.Range("E8").Value = "[Elcio Dec.-Des. Maq. Cad._T2-2B, Elcio Dec.-Des. Maq. Cad._T1-2B]"
